$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1037.625
$ws.Range("I6").Value = 260.2
$ws.Range("K6").Value = 780.5999999999999
$ws.Range("M6").Value = -668.5999999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 874.8570999999999
$ws.Range("I8").Value = 34.8
$ws.Range("J8").Value = 2975
$ws.Range("K8").Value = 104.4
$ws.Range("L8").Value = 8925
$ws.Range("M8").Value = 34.60000000000001
$ws.Range("N8").Value = -9203

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 16669007
$ws.Range("I15").Value = 16669007
$ws.Range("K15").Value = 50007021
$ws.Range("M15").Value = -50006852

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 1314.3334
$ws.Range("I55").Value = 1646.25
$ws.Range("J55").Value = 650.5
$ws.Range("K55").Value = 1646.25
$ws.Range("L55").Value = 650.5
$ws.Range("M55").Value = -1432.25
$ws.Range("N55").Value = -1078.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3859.5454
$ws.Range("I116").Value = 2891.077
$ws.Range("K116").Value = 2891.077
$ws.Range("M116").Value = 550.9229999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5265623
$ws.Range("I132").Value = 6453500.5
$ws.Range("J132").Value = 5022.857
$ws.Range("K132").Value = 19360501.5
$ws.Range("L132").Value = 15068.571
$ws.Range("M132").Value = -19357971.5
$ws.Range("N132").Value = -20128.571

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 594.8
$ws.Range("I135").Value = 694.087
$ws.Range("J135").Value = 268.57144
$ws.Range("K135").Value = 6246.782999999999
$ws.Range("L135").Value = 2417.14296
$ws.Range("M135").Value = -3711.782999999999
$ws.Range("N135").Value = -7487.14296

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3268.1155
$ws.Range("I138").Value = 1700.1818
$ws.Range("J138").Value = 11891.75
$ws.Range("K138").Value = 5100.5454
$ws.Range("L138").Value = 35675.25
$ws.Range("M138").Value = 39.45460000000003
$ws.Range("N138").Value = -45955.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 442041.34
$ws.Range("I141").Value = 1232.2916
$ws.Range("J141").Value = 3968513.8
$ws.Range("K141").Value = 3696.8748
$ws.Range("L141").Value = 11905541.4
$ws.Range("M141").Value = 1483.1252
$ws.Range("N141").Value = -11915901.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1508.2821
$ws.Range("I45").Value = 983.1142599999999
$ws.Range("K45").Value = 983.1142599999999
$ws.Range("M45").Value = -606.1142599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2888.5386
$ws.Range("I61").Value = 974.8125
$ws.Range("J61").Value = 4219.826
$ws.Range("K61").Value = 974.8125
$ws.Range("L61").Value = 4219.826
$ws.Range("M61").Value = -762.8125
$ws.Range("N61").Value = -4643.826

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 744.8333
$ws.Range("I74").Value = 501.73334
$ws.Range("J74").Value = 1150
$ws.Range("K74").Value = 501.73334
$ws.Range("L74").Value = 1150
$ws.Range("M74").Value = 372.26666
$ws.Range("N74").Value = -2898

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 744.8333
$ws.Range("I77").Value = 501.73334
$ws.Range("J77").Value = 1150
$ws.Range("K77").Value = 2508.6667
$ws.Range("L77").Value = 5750
$ws.Range("M77").Value = 1859.3333
$ws.Range("N77").Value = -14486

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 27782126
$ws.Range("I132").Value = 43482904
$ws.Range("J132").Value = 3830.6155
$ws.Range("K132").Value = 130448712
$ws.Range("L132").Value = 11491.8465
$ws.Range("M132").Value = -130446182
$ws.Range("N132").Value = -16551.8465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2888.5386
$ws.Range("I136").Value = 974.8125
$ws.Range("J136").Value = 4219.826
$ws.Range("K136").Value = 2924.4375
$ws.Range("L136").Value = 12659.478
$ws.Range("M136").Value = -374.4375
$ws.Range("N136").Value = -17759.478

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1727.8422
$ws.Range("I99").Value = 1237
$ws.Range("J99").Value = 5900
$ws.Range("K99").Value = 1237
$ws.Range("L99").Value = 5900
$ws.Range("M99").Value = 261
$ws.Range("N99").Value = -8896

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2449.83
$ws.Range("I134").Value = 2224.1707
$ws.Range("K134").Value = 6672.5121
$ws.Range("M134").Value = -4137.5121

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2379.7812
$ws.Range("I31").Value = 1433.641
$ws.Range("J31").Value = 3855.76
$ws.Range("K31").Value = 1433.641
$ws.Range("L31").Value = 3855.76
$ws.Range("M31").Value = -1138.641
$ws.Range("N31").Value = -4445.76

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2379.7812
$ws.Range("I34").Value = 1433.641
$ws.Range("J34").Value = 3855.76
$ws.Range("K34").Value = 1433.641
$ws.Range("L34").Value = 3855.76
$ws.Range("M34").Value = -1231.641
$ws.Range("N34").Value = -4259.76

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4517.3335
$ws.Range("I132").Value = 2837.5
$ws.Range("J132").Value = 5861.2
$ws.Range("K132").Value = 8512.5
$ws.Range("L132").Value = 17583.6
$ws.Range("M132").Value = -5982.5
$ws.Range("N132").Value = -22643.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1778.2545
$ws.Range("I134").Value = 1381.0889
$ws.Range("K134").Value = 4143.2667
$ws.Range("M134").Value = -1608.2667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 7358868
$ws.Range("I139").Value = 10418782
$ws.Range("K139").Value = 31256346
$ws.Range("M139").Value = -31251206

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3606.182
$ws.Range("I70").Value = 3617.9
$ws.Range("J70").Value = 3489
$ws.Range("K70").Value = 3617.9
$ws.Range("L70").Value = 3489
$ws.Range("M70").Value = -3347.9
$ws.Range("N70").Value = -4029

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 3606.182
$ws.Range("I73").Value = 3617.9
$ws.Range("J73").Value = 3489
$ws.Range("K73").Value = 3617.9
$ws.Range("L73").Value = 3489
$ws.Range("M73").Value = -2681.9
$ws.Range("N73").Value = -5361

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2918.4375
$ws.Range("I126").Value = 1911.7646
$ws.Range("J126").Value = 4059.3333
$ws.Range("K126").Value = 5735.293799999999
$ws.Range("L126").Value = 12177.9999
$ws.Range("M126").Value = -3265.293799999999
$ws.Range("N126").Value = -17117.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5620.4546
$ws.Range("I132").Value = 5333.3335
$ws.Range("J132").Value = 5728.125
$ws.Range("K132").Value = 16000.0005
$ws.Range("L132").Value = 17184.375
$ws.Range("M132").Value = -13470.0005
$ws.Range("N132").Value = -22244.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4117.4546
$ws.Range("I132").Value = 2823.5
$ws.Range("J132").Value = 4856.857
$ws.Range("K132").Value = 8470.5
$ws.Range("L132").Value = 14570.571
$ws.Range("M132").Value = -5940.5
$ws.Range("N132").Value = -19630.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 31884
$ws.Range("J133").Value = 31884
$ws.Range("L133").Value = 31884
$ws.Range("N133").Value = -36944

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 846.25
$ws.Range("I81").Value = 650.1
$ws.Range("J81").Value = 1173.1666
$ws.Range("K81").Value = 1300.2
$ws.Range("L81").Value = 2346.3332
$ws.Range("M81").Value = -239.2
$ws.Range("N81").Value = -4468.3332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 846.25
$ws.Range("I84").Value = 650.1
$ws.Range("J84").Value = 1173.1666
$ws.Range("K84").Value = 6501
$ws.Range("L84").Value = 11731.666
$ws.Range("M84").Value = -1197
$ws.Range("N84").Value = -22339.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1400
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1400
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1400
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -4146

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4349739.5
$ws.Range("I126").Value = 1389.1111
$ws.Range("J126").Value = 7145107.5
$ws.Range("K126").Value = 4167.3333
$ws.Range("L126").Value = 21435322.5
$ws.Range("M126").Value = -1697.3333
$ws.Range("N126").Value = -21440262.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5053.6924
$ws.Range("I132").Value = 1659.84
$ws.Range("K132").Value = 4979.52
$ws.Range("M132").Value = -2449.52

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 76871.664
$ws.Range("J135").Value = 76871.664
$ws.Range("L135").Value = 76871.664
$ws.Range("N135").Value = -87011.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1489.2894
$ws.Range("I136").Value = 709.6
$ws.Range("J136").Value = 4413.125
$ws.Range("K136").Value = 2128.8
$ws.Range("L136").Value = 13239.375
$ws.Range("M136").Value = 421.1999999999998
$ws.Range("N136").Value = -18339.375
